$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L: header "WMpass" and value "blackdress19" (same as other password cells)
$ws.Range("L1").Value = "WMpass"
$ws.Range("L2").Value = "blackdress19"

# Update the active selection to match the new last cell
$ws.Range("L2").Select()
